$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Obstical detection BOM" sheet: a new primary component (Raspberry Pi HQ
# camera) replaces the old OPT8241NBN entry; the OPT8241NBN data is kept but
# demoted to the (previously empty) "Secondary Component" column D.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Obstical detection BOM")

$mouserUrl = "https://www.mouser.se/ProductDetail/Texas-Instruments/OPT8241NBN?qs=cGEy3R83DS%2FxFMUAL%252BoBvw%3D%3D"

# Move the existing OPT8241NBN data from column B into column D
$ws6.Range("D1").Value2 = $ws6.Range("B1").Value2
$ws6.Range("D2").Value2 = $ws6.Range("B2").Value2
$ws6.Range("D3").Value2 = $ws6.Range("B3").Value2
$ws6.Range("D4").Value2 = 1
$ws6.Range("D5").Value2 = $ws6.Range("B5").Value2
$ws6.Range("D6").Formula = "=D4*D5"
$ws6.Range("D2").Style = "Hyperlink"
$ws6.Columns.Item(4).ColumnWidth = 41.6080729166667

# Move the Mouser hyperlink from B2 to D2 (clearing the leftover format on B2)
$ws6.Range("B2").Hyperlinks.Delete()
$ws6.Range("B2").ClearFormats()
$ws6.Hyperlinks.Add($ws6.Range("D2"), $mouserUrl) | Out-Null

# Write the new primary component (Raspberry Pi HQ camera) into column B
$ws6.Range("B1").Value2 = "Officiell Raspberry Pi HQ-kamera 12,3 MP ned SONY IMX477R- sensor"
$ws6.Range("B2").Value2 = "r-pi"
$ws6.Range("B3").Value2 = "38*38*18,43)"
$ws6.Range("B5").Value2 = 698

# Hyperlink for the new primary component
$ws6.Hyperlinks.Add($ws6.Range("B2"), "https://www.electrokit.com/raspberry-pi-hq-kamera") | Out-Null

# ---------------------------------------------------------------------------
# "Complete BOM" sheet view state (the F-column formulas pick up the new
# values from "Obstical detection BOM" automatically once recalculated).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Complete BOM")
$ws1.Range("E9").Select()

# ---------------------------------------------------------------------------
# Final view state: "Obstical detection BOM" becomes the active/selected tab.
# ---------------------------------------------------------------------------
$ws6.Activate()
$ws6.Range("B10:B11").Select()
